# Apply the coin-price refresh captured in the commit diff.
# Columns B (Coin) and C (Link) are plain text and can be assigned directly.
# Columns D (Price) and E (Volume) often look numeric (e.g. 1.033, 0.4648)
# but must stay literal text (as in the source inline strings), so we force
# the cell to Text format before assigning, then restore the default "Normal"
# cell style so no stray numeric formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.585.52"
Set-TextValue $ws.Range("E2") "  +2.57%  "

Set-TextValue $ws.Range("D3") "1.852.06"
Set-TextValue $ws.Range("E3") "  +2.09%  "

Set-TextValue $ws.Range("D4") "1.033"
Set-TextValue $ws.Range("E4") "  +2.92%  "

Set-TextValue $ws.Range("D5") "321.90"
Set-TextValue $ws.Range("E5") "  +3.39%  "

Set-TextValue $ws.Range("D6") "1.029"
Set-TextValue $ws.Range("E6") "  +2.54%  "

Set-TextValue $ws.Range("D7") "0.4395"
Set-TextValue $ws.Range("E7") "  +2.42%  "

Set-TextValue $ws.Range("D8") "0.3774"
Set-TextValue $ws.Range("E8") "  +2.30%  "

Set-TextValue $ws.Range("D9") "0.07414"
Set-TextValue $ws.Range("E9") "  +2.45%  "

Set-TextValue $ws.Range("D10") "0.8753"
Set-TextValue $ws.Range("E10") "  +1.61%  "

Set-TextValue $ws.Range("D11") "21.46"
Set-TextValue $ws.Range("E11") "  +1.91%  "

Set-TextValue $ws.Range("D12") "1.863.48"
Set-TextValue $ws.Range("E12") "  -7.74%  "

Set-TextValue $ws.Range("D13") "5.532"
Set-TextValue $ws.Range("E13") "  +2.50%  "

Set-TextValue $ws.Range("D14") "6.694"
Set-TextValue $ws.Range("E14") "  +0.82%  "

Set-TextValue $ws.Range("D15") "0.07221"
Set-TextValue $ws.Range("E15") "  +4.67%  "

Set-TextValue $ws.Range("D16") "82.93"
Set-TextValue $ws.Range("E16") "  +2.79%  "

Set-TextValue $ws.Range("E17") "  +3.18%  "

Set-TextValue $ws.Range("E18") "  +1.17%  "

Set-TextValue $ws.Range("D19") "1.030"
Set-TextValue $ws.Range("E19") "  +2.59%  "

Set-TextValue $ws.Range("E20") "  +1.52%  "

Set-TextValue $ws.Range("D21") "27.599.92"
Set-TextValue $ws.Range("E21") "  +2.43%  "

Set-TextValue $ws.Range("D22") "5.255"
Set-TextValue $ws.Range("E22") "  +1.46%  "

Set-TextValue $ws.Range("E23") "  +2.62%  "

Set-TextValue $ws.Range("D24") "2.073.78"
Set-TextValue $ws.Range("E24") "  -7.78%  "

Set-TextValue $ws.Range("D25") "157.84"
Set-TextValue $ws.Range("E25") "  +2.74%  "

Set-TextValue $ws.Range("D26") "1.923"
Set-TextValue $ws.Range("E26") "  +2.29%  "

Set-TextValue $ws.Range("D27") "18.72"
Set-TextValue $ws.Range("E27") "  +2.71%  "

Set-TextValue $ws.Range("D28") "1.969"
Set-TextValue $ws.Range("E28") "  +5.43%  "

Set-TextValue $ws.Range("D29") "5.258"
Set-TextValue $ws.Range("E29") "  +0.93%  "

Set-TextValue $ws.Range("D30") "117.05"
Set-TextValue $ws.Range("E30") "  +1.82%  "

Set-TextValue $ws.Range("D31") "0.09053"
Set-TextValue $ws.Range("E31") "  +1.24%  "

Set-TextValue $ws.Range("D32") "0.7612"
Set-TextValue $ws.Range("E32") "  +2.47%  "

Set-TextValue $ws.Range("E33") "  +2.66%  "

Set-TextValue $ws.Range("D34") "4.506"
Set-TextValue $ws.Range("E34") "  +1.99%  "

Set-TextValue $ws.Range("D35") "2.885"
Set-TextValue $ws.Range("E35") "  +3.13%  "

Set-TextValue $ws.Range("D36") "1.031"
Set-TextValue $ws.Range("E36") "  +2.28%  "

Set-TextValue $ws.Range("D37") "1.149"
Set-TextValue $ws.Range("E37") "  +2.93%  "

Set-TextValue $ws.Range("D38") "0.01974"
Set-TextValue $ws.Range("E38") "  +2.72%  "

Set-TextValue $ws.Range("D39") "0.05296"
Set-TextValue $ws.Range("E39") "  +1.56%  "

Set-TextValue $ws.Range("D40") "0.5149"
Set-TextValue $ws.Range("E40") "  +1.44%  "

Set-TextValue $ws.Range("D41") "2.807"
Set-TextValue $ws.Range("E41") "  +3.10%  "

Set-TextValue $ws.Range("D42") "0.1674"
Set-TextValue $ws.Range("E42") "  +1.91%  "

Set-TextValue $ws.Range("D43") "6.727"
Set-TextValue $ws.Range("E43") "  +4.74%  "

Set-TextValue $ws.Range("D44") "8.479"
Set-TextValue $ws.Range("E44") "  +2.82%  "

Set-TextValue $ws.Range("D45") "108.74"
Set-TextValue $ws.Range("E45") "  +1.88%  "

Set-TextValue $ws.Range("D46") "10.60"
Set-TextValue $ws.Range("E46") "  +1.81%  "

Set-TextValue $ws.Range("E47") "  +3.07%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.4648"
Set-TextValue $ws.Range("E48") "  +1.63%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.06397"
Set-TextValue $ws.Range("E49") "  +1.74%  "

Set-TextValue $ws.Range("D50") "1.852"
Set-TextValue $ws.Range("E50") "  +2.91%  "

Set-TextValue $ws.Range("D51") "39.12"
Set-TextValue $ws.Range("E51") "  +4.01%  "
